$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textRefs = @("D4", "D5", "D6", "D7", "D9", "D10", "D11", "D13", "D14", "D17", "D19", "D20", "D22", "D23", "D25", "D27", "D29", "D30", "D31", "D32", "D33", "D36", "D38", "D40", "D44", "D48", "D49", "D50", "D51")
foreach ($ref in $textRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "45.236.34"
$ws.Range("E2").Value = "  +2.18%  "
$ws.Range("D3").Value = "2.425.02"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "318.68"
$ws.Range("E5").Value = "  +3.32%  "
$ws.Range("D6").Value = "102.49"
$ws.Range("E6").Value = "  +1.60%  "
$ws.Range("D7").Value = "0.515"
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "0.526"
$ws.Range("E9").Value = "  +5.07%  "
$ws.Range("D10").Value = "35.51"
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("D11").Value = "0.0799"
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("E12").Value = "  -1.81%  "
$ws.Range("D13").Value = "18.23"
$ws.Range("E13").Value = "  -3.12%  "
$ws.Range("D14").Value = "7.07"
$ws.Range("E14").Value = "  +1.92%  "
$ws.Range("D15").Value = "2.805.13"
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("D16").Value = "2.383.48"
$ws.Range("E16").Value = "  -2.26%  "
$ws.Range("D17").Value = "0.845"
$ws.Range("E17").Value = "  +1.34%  "
$ws.Range("D18").Value = "45.161.71"
$ws.Range("E18").Value = "  +2.11%  "
$ws.Range("D19").Value = "12.20"
$ws.Range("E19").Value = "  -0.95%  "
$ws.Range("D20").Value = "6.33"
$ws.Range("E20").Value = "  -1.47%  "
$ws.Range("D21").Value = "0.0₃0921"
$ws.Range("E21").Value = "  +1.63%  "
$ws.Range("D22").Value = "68.93"
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("D23").Value = "244.81"
$ws.Range("E23").Value = "  +1.73%  "
$ws.Range("E24").Value = "  -1.30%  "
$ws.Range("D25").Value = "2.48"
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").Value = "25.76"
$ws.Range("E27").Value = "  +1.96%  "
$ws.Range("E28").Value = "  -2.74%  "
$ws.Range("D29").Value = "9.63"
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("D30").Value = "49.53"
$ws.Range("E30").Value = "  +2.67%  "
$ws.Range("D31").Value = "32.90"
$ws.Range("E31").Value = "  -0.79%  "
$ws.Range("B32").Value = "Celestia"
$ws.Range("C32").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D32").Value = "20.04"
$ws.Range("E32").Value = "  +7.47%  "
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").Value = "0.125"
$ws.Range("E33").Value = "  +7.01%  "
$ws.Range("E34").Value = "  +0.59%  "
$ws.Range("E35").Value = "  +0.25%  "
$ws.Range("D36").Value = "0.0763"
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("E37").Value = "  -1.79%  "
$ws.Range("D38").Value = "4.43"
$ws.Range("E38").Value = "  -1.31%  "
$ws.Range("E39").Value = "  -1.03%  "
$ws.Range("D40").Value = "126.11"
$ws.Range("E40").Value = "  -2.05%  "
$ws.Range("E41").Value = "  +0.47%  "
$ws.Range("E42").Value = "  -3.98%  "
$ws.Range("E43").Value = "  -1.36%  "
$ws.Range("D44").Value = "0.0289"
$ws.Range("E44").Value = "  +0.49%  "
$ws.Range("D45").Value = "1.935.09"
$ws.Range("E45").Value = "  -1.03%  "
$ws.Range("E46").Value = "  -2.74%  "
$ws.Range("E47").Value = "  +2.21%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "9.11"
$ws.Range("E48").Value = "  -3.13%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "1.78"
$ws.Range("E49").Value = "  +7.16%  "
$ws.Range("D50").Value = "76.46"
$ws.Range("E50").Value = "  +3.76%  "
$ws.Range("D51").Value = "4.80"
$ws.Range("E51").Value = "  +5.13%  "
